$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sucrose_c__Day_sp_exchange"
$ws.Range("B2").Value = 0.02497
$ws.Range("C2").Value = -0.296092499999758
$ws.Range("D2").Value = 0.2960925000000199

$ws.Range("A3").Value = "L-Isoleucine__Day_sp_exchange"
$ws.Range("B3").Value = 0.0033
$ws.Range("C3").Value = -0.0033
$ws.Range("D3").Value = 0.003300000000004672

$ws.Range("A4").Value = "L-Leucine__Day_sp_exchange"
$ws.Range("B4").Value = 0.007700000000000001
$ws.Range("C4").Value = -0.007700000000000001
$ws.Range("D4").Value = 0.007700000000000001

$ws.Range("A5").Value = "L-Lysine__Day_sp_exchange"
$ws.Range("B5").Value = 0.000055
$ws.Range("C5").Value = -0.000055
$ws.Range("D5").Value = 0.000055

$ws.Range("A6").Value = "L-Methionine__Day_sp_exchange"
$ws.Range("B6").Value = 0.0001100000000000007
$ws.Range("C6").Value = -0.0001099999999999973
$ws.Range("D6").Value = 0.0001100000000000007

$ws.Range("A7").Value = "L-Phenylalanine__Day_sp_exchange"
$ws.Range("B7").Value = 0.0044
$ws.Range("C7").Value = -0.02793999999999629
$ws.Range("D7").Value = 0.02794000000000039

$ws.Range("A8").Value = "L-Threonine__Day_sp_exchange"
$ws.Range("B8").Value = 0.00473
$ws.Range("C8").Value = -0.01155000000000263
$ws.Range("D8").Value = 0.01155

$ws.Range("A9").Value = "L-Tryptophan__Day_sp_exchange"
$ws.Range("B9").Value = 0.0022
$ws.Range("C9").Value = -0.0022
$ws.Range("D9").Value = 0.0022

$ws.Range("A10").Value = "L-Valine__Day_sp_exchange"
$ws.Range("B10").Value = 0.005500000000000001
$ws.Range("C10").Value = -0.0209
$ws.Range("D10").Value = 0.02090000000001515

$ws.Range("A11").Value = "L-Cysteine__Day_sp_exchange"
$ws.Range("B11").Value = 0.00132
$ws.Range("C11").Value = -0.001320000000052017
$ws.Range("D11").Value = 0.00132

$ws.Range("A12").Value = "L-Glutamine_c__Day_sp_exchange"
$ws.Range("B12").Value = 0.02048199999999901
$ws.Range("C12").Value = -0.1926048928569813
$ws.Range("D12").Value = 0.4810767499997521

$ws.Range("A13").Value = "L-Glutamate_c__Day_sp_exchange"
$ws.Range("B13").Value = 1.416330230157635
$ws.Range("C13").Value = 1.203238055555524
$ws.Range("D13").Value = 7.611845999994399

$ws.Range("A14").Value = "L-Tyrosine__Day_sp_exchange"
$ws.Range("B14").Value = 0.00957
$ws.Range("C14").Value = -0.009570000000006613
$ws.Range("D14").Value = 0.009570000000003095

$ws.Range("A15").Value = "L-Asparagine__Day_sp_exchange"
$ws.Range("B15").Value = 0.02145
$ws.Range("C15").Value = -0.1488162499999055
$ws.Range("D15").Value = 0.1488162500000992

$ws.Range("A16").Value = "L-Serine_c__Day_sp_exchange"
$ws.Range("B16").Value = 0.01155
$ws.Range("C16").Value = -0.2733225000001817
$ws.Range("D16").Value = 0.2733225000001933

$ws.Range("A17").Value = "L-Aspartate_c__Day_sp_exchange"
$ws.Range("B17").Value = 0.2722904206358095
$ws.Range("C17").Value = -4.424533666668475
$ws.Range("D17").Value = 0.6158618888889418

$ws.Range("A18").Value = "Starch_p__Day_sp_exchange"
$ws.Range("B18").Value = 0.1141292777779996
$ws.Range("C18").Value = -0.2856975000007493
$ws.Range("D18").Value = 0.2609575833332883

$ws.Range("A19").Value = "(S)-Malate_c__Day_sp_exchange"
$ws.Range("B19").Value = 1.07362330952184
$ws.Range("C19").Value = 0.9911678333325844
$ws.Range("D19").Value = 9.399137000000328

$ws.Range("A20").Value = "Fumarate__Day_sp_exchange"
$ws.Range("B20").Value = 0.000923999999994571
$ws.Range("C20").Value = -0.03814799999996859
$ws.Range("D20").Value = 0.03814800000000404
